$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column R (year 2021) ------------------------------------------------
# Row 2 (thin separator row under the header) - blank cell, same style as Q2
$ws.Range("Q2").Copy() | Out-Null
$ws.Range("R2").PasteSpecial(-4122) | Out-Null

# Row 3 - year header
$ws.Range("Q3").Copy() | Out-Null
$ws.Range("R3").PasteSpecial(-4122) | Out-Null
$ws.Range("R3").Value = 2021

# Row 4 - headline indicator (bold row), needs its own number format style
$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null
$ws.Range("R4").Value = 18
$ws.Range("R4").NumberFormat = "0.0"

# Row 5
$ws.Range("Q5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null
$ws.Range("R5").Value = 1.7480265877296817

# Row 6
$ws.Range("O6").Copy() | Out-Null
$ws.Range("R6").PasteSpecial(-4122) | Out-Null
$ws.Range("R6").Value = 4.1112601249414027

# Row 7
$ws.Range("O7").Copy() | Out-Null
$ws.Range("R7").PasteSpecial(-4122) | Out-Null
$ws.Range("R7").Value = 1.5225742120245318

# Row 8
$ws.Range("O8").Copy() | Out-Null
$ws.Range("R8").PasteSpecial(-4122) | Out-Null
$ws.Range("R8").Value = 1.2326518235454269

# Row 9
$ws.Range("O9").Copy() | Out-Null
$ws.Range("R9").PasteSpecial(-4122) | Out-Null
$ws.Range("R9").Value = 4.0865392096984241

# Row 10
$ws.Range("O10").Copy() | Out-Null
$ws.Range("R10").PasteSpecial(-4122) | Out-Null
$ws.Range("R10").Value = 0.84876624403485645

# Row 11
$ws.Range("O11").Copy() | Out-Null
$ws.Range("R11").PasteSpecial(-4122) | Out-Null
$ws.Range("R11").Value = 2.1456657699653627

# Row 12
$ws.Range("O12").Copy() | Out-Null
$ws.Range("R12").PasteSpecial(-4122) | Out-Null
$ws.Range("R12").Value = 1.8214779402142154

# Row 13 - bottom (thick-bordered) row
$ws.Range("Q13").Copy() | Out-Null
$ws.Range("R13").PasteSpecial(-4122) | Out-Null
$ws.Range("R13").Value = 0.51989507542472779

# Leave the same selection state the source workbook was saved with.
$ws.Range("R24:R25").Select() | Out-Null
